# Kazuya Mishima Tekken 8 frame data: flip the sign of the "Block" (column E)
# value for a set of moves (tekken_scrap.py data correction).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    2, 18, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33,
    34, 35, 36, 37, 38, 39, 42, 43, 44, 45, 46, 48, 49, 50, 51,
    52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 67, 68,
    69, 70, 71, 72, 75, 76, 77, 79, 80, 81, 83, 85, 87, 88, 90,
    91, 92, 93, 94, 95, 96, 97, 99, 100, 101, 102, 103, 104, 105, 106,
    107, 108, 109, 110, 111
)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $old = $cell.Value2
    $cell.Value = (0 - $old)
}
